$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.6 = 9641.3 pesos`n✅ 9641.3 pesos = 2.59 = 946.61 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $conversionText

# --- tasas: update the rate/figures used by the "transfi" block ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 384.7
$wsTasas.Range("O10").Value = 3709.01
$wsTasas.Range("N12").Value = 3728.45
$wsTasas.Range("O12").Value = 366.07
